$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A51").Value = 50
$ws.Range("B51").Value = 85
$ws.Range("C51").Value = 1
$ws.Range("D51").Value = 10
$ws.Range("E51").Value = 22
$ws.Range("F51").Value = 96
$ws.Range("G51").Value = 118
